# Applies the December inventory cashflow update to Sheet1.
# Updates J, P, Q, R, S, V columns for the rows affected by the
# reconciliation recalculation described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 4576
$ws.Range("P2").Value = 3853
$ws.Range("Q2").Value = 510.29
$ws.Range("R2").Value = 5108
$ws.Range("V2").Value = 4576

# Row 3
$ws.Range("P3").Value = 26962
$ws.Range("Q3").Value = 71.26000000000001

# Row 4
$ws.Range("P4").Value = 21152
$ws.Range("R4").Value = 259
$ws.Range("S4").Value = 27

# Row 5
$ws.Range("P5").Value = 144550

# Row 6
$ws.Range("J6").Value = 262
$ws.Range("P6").Value = 38454
$ws.Range("Q6").Value = 32.3
$ws.Range("S6").Value = 46
$ws.Range("V6").Value = 262

# Row 7
$ws.Range("J7").Value = 790
$ws.Range("P7").Value = 12903
$ws.Range("Q7").Value = 92.45999999999999
$ws.Range("R7").Value = 937
$ws.Range("S7").Value = 149
$ws.Range("V7").Value = 790

# Row 8
$ws.Range("P8").Value = 152271
$ws.Range("Q8").Value = 35.74

# Row 10
$ws.Range("P10").Value = 3853

# Row 11
$ws.Range("P11").Value = 3853

# Row 12
$ws.Range("Q12").Value = 61.91

# Row 13
$ws.Range("P13").Value = 140624
$ws.Range("Q13").Value = 20.08

# Row 15
$ws.Range("P15").Value = 103756
$ws.Range("Q15").Value = 42.53

# Row 16
$ws.Range("J16").Value = 9838
$ws.Range("P16").Value = 652231
$ws.Range("Q16").Value = 1117.31
$ws.Range("R16").Value = 10913
$ws.Range("S16").Value = 1100
$ws.Range("V16").Value = 9838
